$d = $word.ActiveDocument

# Locate the end of "...so I threw it out and restarted. " (right before
# "In addition, I have some basic file reading mechanisms...") and insert
# the new sentence about smart pointers / move semantics there.
$anchor = $d.Content
$anchor.Find.Execute("so I threw it out and restarted. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)
$anchor.InsertAfter("One major thing I changed was I began to use smart pointers and move semantics. ")

# Force the newly inserted text to live in its own run (distinct from the
# runs before/after it) even though its formatting ends up identical, by
# toggling a character property on and back off.
$newRun = $d.Content
$newRun.Find.Execute("One major thing I changed was I began to use smart pointers and move semantics. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newRun.Bold = 1
$newRun.Bold = 0

Write-Output "Inserted new sentence into the 4/10 update paragraph."
